$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Old Data")
$ws1.Range("E2").Value = "S38"
$ws1.Range("E3").Value = "S219"
$ws1.Range("E5").Value = "S36"
$ws1.Range("E7").Value = "S44"
$ws1.Range("E8").Value = "S36"
$ws1.Range("I9").Value = "S00000000"
$ws1.Range("E15").Value = "S22"
$ws1.Range("E16").Value = "S18"
$ws1.Range("E19").Value = "S4"
$ws1.Range("I19").Value = "S00000000"
$ws1.Range("E20").Value = "S72"
$ws1.Range("E21").Value = "S40"
$ws1.Range("E22").Value = "S2"
$ws1.Range("F22").Value = "S177808"
$ws1.Range("D23").Value = "S100"
$ws1.Range("F23").Value = "S158130"
$ws1.Range("E25").Value = "S137"
$ws1.Range("F25").Value = "S159297"
$ws1.Range("E31").Value = "S47"
$ws1.Range("E32").Value = "S14"
$ws1.Range("E33").Value = "S1"
$ws1.Range("E36").Value = "S15"
$ws1.Range("E37").Value = "S4"
$ws1.Range("E38").Value = "S36"
$ws1.Range("E39").Value = "S25"
$ws1.Range("F39").Value = "S213882"
$ws1.Range("E41").Value = "S20"
$ws1.Range("E43").Value = "S33"
$ws1.Range("E47").Value = "S19"
$ws1.Range("E50").Value = "S11"
$ws1.Range("D51").Value = "S20"
$ws1.Range("E51").Value = "S19"
$ws1.Range("F51").Value = "S197446"
$ws1.Range("E54").Value = "S29"
$ws1.Range("E57").Value = "S37"
$ws1.Range("F57").Value = "S211952"
$ws1.Range("E58").Value = "S7"
$ws1.Range("E59").Value = "S28"
$ws1.Range("E60").Value = "S7"
$ws1.Range("E61").Value = "S5"
$ws1.Range("E65").Value = "S20"
$ws1.Range("E66").Value = "S30"
$ws1.Range("E71").Value = "S89"
$ws1.Range("E73").Value = "S26"
$ws1.Range("E74").Value = "S3"
$ws1.Range("E76").Value = "S35"
$ws1.Range("E78").Value = "S45"
$ws1.Range("F78").Value = "S224644"
$ws1.Range("E79").Value = "S30"
$ws1.Range("E80").Value = "S6"
$ws1.Range("E91").Value = "S15"
$ws1.Range("F93").Value = "S137515"
$ws1.Range("E94").Value = "S98"
$ws1.Range("E95").Value = "S18"
$ws1.Range("D96").Value = "S180"
$ws1.Range("E96").Value = "S146"
$ws1.Range("F96").Value = "S138724"
$ws1.Range("E97").Value = "S34"
$ws1.Range("E99").Value = "S14"
$ws1.Range("E101").Value = "S8"
$ws1.Range("E102").Value = "S60"
$ws1.Range("E103").Value = "S162"
$ws1.Range("E104").Value = "S2"
$ws1.Range("D106").Value = "S180"
$ws1.Range("E106").Value = "S118"
$ws1.Range("F106").Value = "S150270"
$ws1.Range("E107").Value = "S9"
$ws1.Range("F108").Value = "S187764"
$ws1.Range("E109").Value = "S22"
$ws1.Range("F110").Value = "S188858"
$ws1.Range("E111").Value = "S8"
$ws1.Range("E114").Value = "S36"
$ws1.Range("E120").Value = "S13"
$ws1.Range("F121").Value = "S182026"
$ws1.Range("F123").Value = "S172127"
$ws1.Range("F126").Value = "S198650"
$ws1.Range("E128").Value = "S24"

$ws2 = $wb.Worksheets.Item("New Data")
$ws2.Range("E2").Value = "S38"
$ws2.Range("E3").Value = "S219"
$ws2.Range("E5").Value = "S36"
$ws2.Range("E7").Value = "S44"
$ws2.Range("E8").Value = "S36"
$ws2.Range("I9").Value = "S00000000"
$ws2.Range("E15").Value = "S22"
$ws2.Range("E16").Value = "S18"
$ws2.Range("E19").Value = "S4"
$ws2.Range("I19").Value = "S00000000"
$ws2.Range("E20").Value = "S72"
$ws2.Range("E21").Value = "S40"
$ws2.Range("E22").Value = "S2"
$ws2.Range("F22").Value = "S177808"
$ws2.Range("D23").Value = "S100"
$ws2.Range("E23").Value = "S14"
$ws2.Range("F23").Value = "S158130"
$ws2.Range("E25").Value = "S137"
$ws2.Range("F25").Value = "S159297"
$ws2.Range("E31").Value = "S47"
$ws2.Range("E32").Value = "S14"
$ws2.Range("E33").Value = "S1"
$ws2.Range("E37").Value = "S4"
$ws2.Range("E38").Value = "S36"
$ws2.Range("E39").Value = "S25"
$ws2.Range("F39").Value = "S213882"
$ws2.Range("E41").Value = "S20"
$ws2.Range("E43").Value = "S33"
$ws2.Range("E47").Value = "S19"
$ws2.Range("E50").Value = "S11"
$ws2.Range("D51").Value = "S20"
$ws2.Range("E51").Value = "S19"
$ws2.Range("F51").Value = "S197446"
$ws2.Range("E54").Value = "S29"
$ws2.Range("E57").Value = "S37"
$ws2.Range("F57").Value = "S211952"
$ws2.Range("E58").Value = "S7"
$ws2.Range("E59").Value = "S28"
$ws2.Range("E60").Value = "S7"
$ws2.Range("E65").Value = "S20"
$ws2.Range("E66").Value = "S30"
$ws2.Range("E71").Value = "S89"
$ws2.Range("E73").Value = "S26"
$ws2.Range("E74").Value = "S3"
$ws2.Range("E76").Value = "S35"
$ws2.Range("E78").Value = "S45"
$ws2.Range("F78").Value = "S224644"
$ws2.Range("E79").Value = "S30"
$ws2.Range("E80").Value = "S6"
$ws2.Range("E91").Value = "S15"
$ws2.Range("D96").Value = "S180"
$ws2.Range("E96").Value = "S146"
$ws2.Range("F96").Value = "S138724"
$ws2.Range("E97").Value = "S34"
$ws2.Range("E101").Value = "S8"
$ws2.Range("E102").Value = "S60"
$ws2.Range("E103").Value = "S162"
$ws2.Range("E104").Value = "S2"
$ws2.Range("D106").Value = "S180"
$ws2.Range("E106").Value = "S118"
$ws2.Range("F106").Value = "S150270"
$ws2.Range("F108").Value = "S187764"
$ws2.Range("E109").Value = "S22"
$ws2.Range("E111").Value = "S8"
$ws2.Range("E114").Value = "S36"
$ws2.Range("E120").Value = "S13"
$ws2.Range("F126").Value = "S198650"
$ws2.Range("E128").Value = "S24"
